$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore C10 to 1 (was 18)
$ws.Range("C10").Value = 1
